$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = 'Datos actualizados a 21 de Marzo de 2020 a las 14:16'

# --- Country name re-shuffle (column A) ---
$ws.Range("A43").Value = 'Arabia Saudita'
$ws.Range("A44").Value = 'Eslovenia'
$ws.Range("A45").Value = 'Rumania'
$ws.Range("A56").Value = 'Libano'
$ws.Range("A57").Value = 'Mexico'
$ws.Range("A58").Value = 'Panama'
$ws.Range("A80").Value = 'Republica de Macedonia'
$ws.Range("A81").Value = 'Jordania'
$ws.Range("A82").Value = 'Brunei'
$ws.Range("A85").Value = 'Bielorrusia'
$ws.Range("A86").Value = 'Republica de Chipre'
$ws.Range("A87").Value = 'Malta'
$ws.Range("A88").Value = 'Republica Dominicana'
$ws.Range("A89").Value = 'Lituania'
$ws.Range("A98").Value = 'Camboya'
$ws.Range("A99").Value = 'Guadalupe'
$ws.Range("A114").Value = 'Ghana'
$ws.Range("A115").Value = 'Bolivia'
$ws.Range("A120").Value = 'Guayana Francesa'
$ws.Range("A121").Value = 'Polinesia Francesa'
$ws.Range("A122").Value = 'Guam'
$ws.Range("A123").Value = 'Kirguistan'
$ws.Range("A124").Value = 'Puerto Rico'

# --- Updated case-count figures (columns B-H) ---
$ws.Range("B9").Value = 19777
$ws.Range("C9").Value = 394
$ws.Range("E9").Value = 19354
$ws.Range("B14").Value = 3631
$ws.Range("C14").Value = 637
$ws.Range("E14").Value = 3493
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 136
$ws.Range("B18").Value = 1746
$ws.Range("C18").Value = 107
$ws.Range("E18").Value = 1710
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 20
$ws.Range("B25").Value = 986
$ws.Range("C25").Value = 16
$ws.Range("E25").Value = 973
$ws.Range("F26").Value = 7
$ws.Range("B32").Value = 666
$ws.Range("C32").Value = 165
$ws.Range("E32").Value = 650
$ws.Range("E33").Value = 510
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 1
$ws.Range("B43").Value = 392
$ws.Range("C43").Value = 48
$ws.Range("D43").Value = 16
$ws.Range("E43").Value = 376
$ws.Range("F43").Value = 0
$ws.Range("B44").Value = 383
$ws.Range("C44").Value = 42
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 382
$ws.Range("F44").Value = 9
$ws.Range("H44").Value = 1
$ws.Range("B45").Value = 367
$ws.Range("C45").Value = 59
$ws.Range("D45").Value = 52
$ws.Range("E45").Value = 315
$ws.Range("F45").Value = 14
$ws.Range("H45").Value = 0
$ws.Range("B56").Value = 206
$ws.Range("C56").Value = 29
$ws.Range("D56").Value = 8
$ws.Range("E56").Value = 194
$ws.Range("F56").Value = 4
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 4
$ws.Range("B57").Value = 203
$ws.Range("C57").Value = 39
$ws.Range("D57").Value = 4
$ws.Range("E57").Value = 197
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 2
$ws.Range("B58").Value = 200
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 1
$ws.Range("E58").Value = 198
$ws.Range("F58").Value = 7
$ws.Range("H58").Value = 1
$ws.Range("B80").Value = 85
$ws.Range("C80").Value = 9
$ws.Range("E80").Value = 84
$ws.Range("F80").Value = 1
$ws.Range("B81").Value = 84
$ws.Range("C81").Value = 0
$ws.Range("E81").Value = 83
$ws.Range("F81").Value = 0
$ws.Range("B82").Value = 83
$ws.Range("C82").Value = 5
$ws.Range("E82").Value = 82
$ws.Range("F82").Value = 2
$ws.Range("B85").Value = 76
$ws.Range("C85").Value = 7
$ws.Range("D85").Value = 15
$ws.Range("E85").Value = 61
$ws.Range("F85").Value = 0
$ws.Range("B86").Value = 75
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 75
$ws.Range("B87").Value = 73
$ws.Range("C87").Value = 9
$ws.Range("D87").Value = 2
$ws.Range("E87").Value = 71
$ws.Range("F87").Value = 1
$ws.Range("H87").Value = 0
$ws.Range("B88").Value = 72
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 70
$ws.Range("F88").Value = 0
$ws.Range("H88").Value = 2
$ws.Range("D89").Value = 1
$ws.Range("E89").Value = 67
$ws.Range("F89").Value = 1
$ws.Range("H89").Value = 1
$ws.Range("D98").Value = 1
$ws.Range("F98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("F99").Value = 4
$ws.Range("H99").Value = 1
$ws.Range("C120").Value = 0
$ws.Range("C121").Value = 4
$ws.Range("C122").Value = 1
$ws.Range("C123").Value = 8
$ws.Range("C124").Value = 0
